$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row (19) below the existing last row (18), mirroring row 18's
# layout: a label in column A and a hyperlinked "http://google.com" in column B.
$ws.Range("A19").Value = "iegcViolMsgsFetchUrl"

$ws.Hyperlinks.Add($ws.Range("B19"), "http://google.com", "") | Out-Null

# Hyperlinks.Add re-styles the cell on its own; reapply the same "Hyperlink"
# cell style already used by the other links in column B (e.g. B18) so the
# new cell matches them exactly.
$ws.Range("B19").Style = "Hyperlink"

# Match the active selection left behind by the recorded session.
$ws.Range("K21").Select() | Out-Null
